$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p081v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p081v_1</id>", 2)
$d.Content.Find.Execute("<id>p081v_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p081v_2</id>", 2)
$d.Content.Find.Execute("<id>p081v_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p081v_3</id>", 2)
$d.Content.Find.Execute("<id>p081v_a4</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p081v_4</id>", 2)
